# Sketch out the EPFParser interface on the "EPFImporter" notes sheet,
# between the existing EPFIngester block and the EPFDbUtil block.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EPFImporter")

# Insert two new rows above what is currently row 19, pushing the
# parseTableName/parseColumnsAndTypes/.../nextRecord() rows and the whole
# EPFDbUtil block below them down by two rows.
$ws.Range("19:20").Insert()

# Use the first inserted row as a new section header, mirroring how
# "EPFIngester" (A16) and "EPFDbUtil" introduce their own method lists.
$ws.Range("A20").Value = "EPFParser"

# Leave the sheet scrolled to/selected on this new section, matching the
# saved view state.
$ws.Activate()
$ws.Range("B19").Select()
$excel.ActiveWindow.ScrollRow = 10
